$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record the existing hyperlinks (address, target, display text) before the
# row shift, since the hyperlink ranges don't automatically follow a row
# delete.
$existingLinks = @()
foreach ($hl in $ws.Hyperlinks) {
    $existingLinks += ,@($hl.Range.Row, $hl.Range.Column, $hl.Address, $hl.TextToDisplay)
}

# Remove the first review row (row 2), shifting the rest of the rows up.
$ws.Range("A2:H2").EntireRow.Delete()

# Drop the now-stale hyperlinks and re-create them one row higher so they
# keep pointing at the same (shifted) cells.
$ws.Hyperlinks.Delete()
foreach ($link in $existingLinks) {
    $newRow = $link[0] - 1
    $col = $link[1]
    $target = $link[2]
    $display = $link[3]
    $cell = $ws.Cells.Item($newRow, $col)
    $ws.Hyperlinks.Add($cell, $target, "", "", $display) | Out-Null
}
